$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.653.09'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '2.235.12'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = "'305.48"
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').Value = "'94.86"
$ws.Range('E6').Value = '  -1.05%  '
$ws.Range('D7').Value = "'0.572"
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = "'0.518"
$ws.Range('E9').Value = '  -1.73%  '
$ws.Range('D10').Value = "'34.85"
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('E11').Value = '  -1.96%  '
$ws.Range('D12').Value = "'7.17"
$ws.Range('E12').Value = '  -1.27%  '
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').Value = '2.579.54'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').Value = '2.245.78'
$ws.Range('E15').Value = '  -4.90%  '
$ws.Range('D16').Value = "'0.835"
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = "'13.52"
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '44.478.11'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('D19').Value = '0.0₃0942'
$ws.Range('E19').Value = '  -3.51%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').Value = "'11.92"
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = "'6.22"
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').Value = "'65.25"
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('D23').Value = "'238.75"
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('D24').Value = "'2.94"
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('E25').Value = '  -1.85%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +4.77%  '
$ws.Range('D28').Value = "'9.79"
$ws.Range('E28').Value = '  -1.59%  '
$ws.Range('D29').Value = "'37.55"
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = "'150.48"
$ws.Range('E32').Value = '  -1.21%  '
$ws.Range('E33').Value = '  -1.79%  '
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('E35').Value = '  -8.85%  '
$ws.Range('D36').Value = "'0.118"
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('D37').Value = "'0.108"
$ws.Range('E37').Value = '  -2.24%  '
$ws.Range('D38').Value = "'1.83"
$ws.Range('E38').Value = '  +3.56%  '
$ws.Range('D39').Value = "'15.07"
$ws.Range('E39').Value = '  +3.31%  '
$ws.Range('D40').Value = "'3.36"
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').Value = "'3.75"
$ws.Range('E42').Value = '  -3.42%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').Value = '1.827.53'
$ws.Range('E44').Value = '  +5.17%  '
$ws.Range('D45').Value = "'1.72"
$ws.Range('E45').Value = '  +10.70%  '
$ws.Range('D46').Value = "'79.91"
$ws.Range('E46').Value = '  -4.22%  '
$ws.Range('D47').Value = "'0.188"
$ws.Range('E47').Value = '  -1.27%  '
$ws.Range('D48').Value = "'98.35"
$ws.Range('E48').Value = '  -2.04%  '
$ws.Range('D49').Value = "'4.87"
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('D50').Value = "'68.86"
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('D51').Value = "'54.09"
$ws.Range('E51').Value = '  -1.26%  '
